# Applies the commit "Updated cryptos list on Fri May  5 09:58:49 UTC 2023 with GitHub Actions"
# Refreshes the Price (D) and Volume(1h) (E) columns for rows 2-51 with the latest crypto feed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.147.22'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.903.20'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.59'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4596'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3888'
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07872'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9892'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.90'
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.883.03'
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.766'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.046'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07004'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.93'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.07'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.161.84'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.316'
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.12'
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.102'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.01'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.44'
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.893'
$ws.Range("E27").Value = '  -3.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '118.42'
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.866'
$ws.Range("E29").Value = '  -6.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09337'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8923'
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.240'
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.320'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.144'
$ws.Range("E34").Value = '  -3.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05778'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.168'
$ws.Range("E36").Value = '  -3.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02088'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.001'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.666'
$ws.Range("E39").Value = '  -3.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5676'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1795'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.708'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.85'
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5347'
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.200'
$ws.Range("E45").Value = '  -1.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07014'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.845'
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.09'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.2937'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("E51").Value = '  -2.53%  '
